$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H32").Value = 1349
$ws.Range("I32").Value = 601
$ws.Range("J32").Value = 1455.8572
$ws.Range("K32").Value = 601
$ws.Range("L32").Value = 1455.8572
$ws.Range("M32").Value = -275
$ws.Range("N32").Value = -2107.8572

$ws.Range("H34").Value = 4268.6
$ws.Range("I34").Value = 1823.25
$ws.Range("J34").Value = 14050
$ws.Range("K34").Value = 1823.25
$ws.Range("L34").Value = 14050
$ws.Range("M34").Value = -1620.25
$ws.Range("N34").Value = -14456

$ws.Range("H36").Value = 4268.6
$ws.Range("I36").Value = 1823.25
$ws.Range("J36").Value = 14050
$ws.Range("K36").Value = 1823.25
$ws.Range("L36").Value = 14050
$ws.Range("M36").Value = -1108.25
$ws.Range("N36").Value = -15480

$ws.Range("H47").Value = 14066.667
$ws.Range("J47").Value = 14066.667
$ws.Range("L47").Value = 14066.667
$ws.Range("N47").Value = -16010.667

$ws.Range("H54").Value = 14000
$ws.Range("I54").Value = 9000
$ws.Range("J54").Value = 15000
$ws.Range("K54").Value = 9000
$ws.Range("L54").Value = 15000
$ws.Range("M54").Value = -8514
$ws.Range("N54").Value = -15972

$ws.Range("H100").Value = 41668544
$ws.Range("I100").Value = 55556892
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 55556892
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -55556351
$ws.Range("N100").Value = -4582

$ws.Range("H106").Value = 333338080
$ws.Range("I106").Value = 111117450
$ws.Range("K106").Value = 111117450
$ws.Range("M106").Value = -111116819

$ws.Range("H107").Value = 20838922
$ws.Range("I107").Value = 25001566
$ws.Range("J107").Value = 25703
$ws.Range("K107").Value = 25001566
$ws.Range("L107").Value = 25703
$ws.Range("M107").Value = -24999646
$ws.Range("N107").Value = -29543

$ws.Range("H111").Value = 102860.1
$ws.Range("I111").Value = 3750
$ws.Range("J111").Value = 168933.5
$ws.Range("K111").Value = 11250
$ws.Range("L111").Value = 506800.5
$ws.Range("M111").Value = -8183
$ws.Range("N111").Value = -512934.5

$ws.Range("H113").Value = 15328.333
$ws.Range("I113").Value = 15328.333
$ws.Range("K113").Value = 15328.333
$ws.Range("M113").Value = -12074.333

$ws.Range("H137").Value = 1893.6
$ws.Range("I137").Value = 1888.625
$ws.Range("J137").Value = 1899.2858
$ws.Range("K137").Value = 5665.875
$ws.Range("L137").Value = 5697.857400000001
$ws.Range("M137").Value = -3115.875
$ws.Range("N137").Value = -10797.8574

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 3200
$ws.Range("I45").Value = 3625
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 3625
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -3248
$ws.Range("N45").Value = -2254

$ws.Range("H74").Value = 4251.6216
$ws.Range("I74").Value = 7812.3335
$ws.Range("K74").Value = 7812.3335
$ws.Range("M74").Value = -6938.3335

$ws.Range("H77").Value = 4251.6216
$ws.Range("I77").Value = 7812.3335
$ws.Range("K77").Value = 39061.6675
$ws.Range("M77").Value = -34693.6675

$ws.Range("H122").Value = 1833067.6
$ws.Range("I122").Value = 2138262.2
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 6414786.600000001
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -6412336.600000001
$ws.Range("N122").Value = -10600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H99").Value = 58825588
$ws.Range("I99").Value = 90910820
$ws.Range("J99").Value = 2666.6667
$ws.Range("K99").Value = 90910820
$ws.Range("L99").Value = 2666.6667
$ws.Range("M99").Value = -90909322
$ws.Range("N99").Value = -5662.6667

$ws.Range("H107").Value = 1866.5
$ws.Range("I107").Value = 1839.8
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1839.8
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 80.20000000000005
$ws.Range("N107").Value = -5840

$ws.Range("H134").Value = 3592.3389
$ws.Range("I134").Value = 3917.476
$ws.Range("K134").Value = 11752.428
$ws.Range("M134").Value = -9217.428

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H99").Value = 6252955.5
$ws.Range("I99").Value = 1599.4667
$ws.Range("J99").Value = 25007022
$ws.Range("K99").Value = 1599.4667
$ws.Range("L99").Value = 25007022
$ws.Range("M99").Value = -101.4666999999999
$ws.Range("N99").Value = -25010018

$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524

$ws.Range("H107").Value = 719.55554
$ws.Range("I107").Value = 622.2353000000001
$ws.Range("K107").Value = 622.2353000000001
$ws.Range("M107").Value = 1297.7647

$ws.Range("H109").Value = 49642.5
$ws.Range("J109").Value = 49642.5
$ws.Range("L109").Value = 49642.5
$ws.Range("N109").Value = -51722.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 6252955.5
$ws.Range("I126").Value = 1599.4667
$ws.Range("J126").Value = 25007022
$ws.Range("K126").Value = 4798.4001
$ws.Range("L126").Value = 75021066
$ws.Range("M126").Value = -2328.4001
$ws.Range("N126").Value = -75026006

$ws.Range("H132").Value = 2103.2188
$ws.Range("I132").Value = 1677.625
$ws.Range("J132").Value = 3380
$ws.Range("K132").Value = 5032.875
$ws.Range("L132").Value = 10140
$ws.Range("M132").Value = -2502.875
$ws.Range("N132").Value = -15200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H122").Value = 5895137.5
$ws.Range("I122").Value = 7203501.5
$ws.Range("K122").Value = 21610504.5
$ws.Range("M122").Value = -21608054.5

$ws.Range("H126").Value = 4714.1177
$ws.Range("I126").Value = 6105.9565
$ws.Range("J126").Value = 1803.909
$ws.Range("K126").Value = 18317.8695
$ws.Range("L126").Value = 5411.727000000001
$ws.Range("M126").Value = -15847.8695
$ws.Range("N126").Value = -10351.727

$ws.Range("H132").Value = 4695.9165
$ws.Range("I132").Value = 9713.5
$ws.Range("J132").Value = 3692.4
$ws.Range("K132").Value = 29140.5
$ws.Range("L132").Value = 11077.2
$ws.Range("M132").Value = -26610.5
$ws.Range("N132").Value = -16137.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 85815.664
$ws.Range("I7").Value = 127123.5
$ws.Range("J7").Value = 3200
$ws.Range("K7").Value = 127123.5
$ws.Range("L7").Value = 3200
$ws.Range("M7").Value = -127011.5
$ws.Range("N7").Value = -3424

$ws.Range("H40").Value = 500009000
$ws.Range("I40").Value = 500009000
$ws.Range("K40").Value = 500009000
$ws.Range("M40").Value = -500008864

$ws.Range("H46").Value = 33334636
$ws.Range("I46").Value = 66667650
$ws.Range("J46").Value = 1624.4
$ws.Range("K46").Value = 66667650
$ws.Range("L46").Value = 1624.4
$ws.Range("M46").Value = -66667462
$ws.Range("N46").Value = -2000.4

$ws.Range("H61").Value = 2947
$ws.Range("I61").Value = 2169.25
$ws.Range("J61").Value = 4502.5
$ws.Range("K61").Value = 2169.25
$ws.Range("L61").Value = 4502.5
$ws.Range("M61").Value = -1967.25
$ws.Range("N61").Value = -4906.5

$ws.Range("H113").Value = 2947
$ws.Range("I113").Value = 2169.25
$ws.Range("J113").Value = 4502.5
$ws.Range("K113").Value = 2169.25
$ws.Range("L113").Value = 4502.5
$ws.Range("M113").Value = 0.75
$ws.Range("N113").Value = -8842.5

$ws.Range("H122").Value = 11635367
$ws.Range("I122").Value = 14288914
$ws.Range("J122").Value = 5001500
$ws.Range("K122").Value = 42866742
$ws.Range("L122").Value = 15004500
$ws.Range("M122").Value = -42864292
$ws.Range("N122").Value = -15009400

$ws.Range("H126").Value = 85815.664
$ws.Range("I126").Value = 127123.5
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 381370.5
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -378900.5
$ws.Range("N126").Value = -14540

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 2147.3684
$ws.Range("I122").Value = 1126.6666
$ws.Range("K122").Value = 3379.9998
$ws.Range("M122").Value = -929.9998000000001

$ws.Range("H126").Value = 1398.7273
$ws.Range("I126").Value = 1172.625
$ws.Range("J126").Value = 2001.6666
$ws.Range("K126").Value = 3517.875
$ws.Range("L126").Value = 6004.9998
$ws.Range("M126").Value = -1047.875
$ws.Range("N126").Value = -10944.9998

$ws.Range("H132").Value = 3772.5715
$ws.Range("I132").Value = 7536
$ws.Range("J132").Value = 2746.182
$ws.Range("K132").Value = 22608
$ws.Range("L132").Value = 8238.545999999998
$ws.Range("M132").Value = -20078
$ws.Range("N132").Value = -13298.546
